$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.NumberFormat = "General"
    $cell.Style = "Normal"
}

# --- Rows 30/31 swap: ImmutableX <-> Mantle, with updated price/volume ---
$ws.Range("B30").Value = "Mantle"
$ws.Range("C30").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue $ws.Range("D30") "1.15"
$ws.Range("E30").Value = "  +3.97%  "

$ws.Range("B31").Value = "ImmutableX"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue $ws.Range("D31") "1.87"
$ws.Range("E31").Value = "  -0.40%  "

# --- Column D (Price) updates ---
Set-TextValue $ws.Range("D2") "60.874.28"
Set-TextValue $ws.Range("D3") "2.971.64"
Set-TextValue $ws.Range("D5") "533.89"
Set-TextValue $ws.Range("D6") "133.14"
Set-TextValue $ws.Range("D8") "2.965.38"
Set-TextValue $ws.Range("D9") "0.492"
Set-TextValue $ws.Range("D10") "0.147"
Set-TextValue $ws.Range("D11") "6.07"
Set-TextValue $ws.Range("D12") "0.442"
Set-TextValue $ws.Range("D13") "0.0000219"
Set-TextValue $ws.Range("D14") "33.78"
Set-TextValue $ws.Range("D15") "3.474.08"
Set-TextValue $ws.Range("D17") "61.035.32"
Set-TextValue $ws.Range("D18") "2.996.93"
Set-TextValue $ws.Range("D19") "6.54"
Set-TextValue $ws.Range("D20") "462.31"
Set-TextValue $ws.Range("D21") "13.12"
Set-TextValue $ws.Range("D22") "0.668"
Set-TextValue $ws.Range("D23") "6.88"
Set-TextValue $ws.Range("D24") "79.24"
Set-TextValue $ws.Range("D25") "11.91"
Set-TextValue $ws.Range("D27") "2.66"
Set-TextValue $ws.Range("D28") "7.78"
Set-TextValue $ws.Range("D32") "25.31"
Set-TextValue $ws.Range("D33") "55.26"
Set-TextValue $ws.Range("D34") "5.39"
Set-TextValue $ws.Range("D35") "2.25"
Set-TextValue $ws.Range("D36") "5.81"
Set-TextValue $ws.Range("D37") "449.52"
Set-TextValue $ws.Range("D38") "3.185.92"
Set-TextValue $ws.Range("D39") "0.0780"
Set-TextValue $ws.Range("D40") "0.0380"
Set-TextValue $ws.Range("D41") "0.117"
Set-TextValue $ws.Range("D42") "8.06"
Set-TextValue $ws.Range("D43") "27.21"
Set-TextValue $ws.Range("D44") "2.45"
Set-TextValue $ws.Range("D46") "0.241"
Set-TextValue $ws.Range("D47") "119.42"
Set-TextValue $ws.Range("D48") "1.98"
Set-TextValue $ws.Range("D50") "0.0₃0492"

# --- Column E (Volume 1h) updates ---
$ws.Range("E2").Value = "  -3.22%  "
$ws.Range("E3").Value = "  -3.06%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("E5").Value = "  -0.63%  "
$ws.Range("E6").Value = "  -0.48%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("E8").Value = "  -3.02%  "
$ws.Range("E9").Value = "  -0.67%  "
$ws.Range("E10").Value = "  -4.63%  "
$ws.Range("E11").Value = "  +0.30%  "
$ws.Range("E12").Value = "  -2.23%  "
$ws.Range("E13").Value = "  -2.07%  "
$ws.Range("E14").Value = "  -1.36%  "
$ws.Range("E15").Value = "  -2.39%  "
$ws.Range("E16").Value = "  -1.01%  "
$ws.Range("E17").Value = "  -2.98%  "
$ws.Range("E18").Value = "  -2.10%  "
$ws.Range("E19").Value = "  -1.43%  "
$ws.Range("E20").Value = "  -4.21%  "
$ws.Range("E21").Value = "  -1.59%  "
$ws.Range("E22").Value = "  -3.91%  "
$ws.Range("E23").Value = "  -3.20%  "
$ws.Range("E24").Value = "  +0.05%  "
$ws.Range("E25").Value = "  -1.56%  "
$ws.Range("E26").Value = "  -0.11%  "
$ws.Range("E27").Value = "  -1.46%  "
$ws.Range("E28").Value = "  -4.43%  "
$ws.Range("E29").Value = "  +0.31%  "
$ws.Range("E32").Value = "  -2.62%  "
$ws.Range("E33").Value = "  -3.19%  "
$ws.Range("E34").Value = "  +0.89%  "
$ws.Range("E35").Value = "  -4.72%  "
$ws.Range("E36").Value = "  -3.49%  "
$ws.Range("E37").Value = "  -7.62%  "
$ws.Range("E38").Value = "  +1.66%  "
$ws.Range("E39").Value = "  -1.70%  "
$ws.Range("E40").Value = "  -3.71%  "
$ws.Range("E41").Value = "  +1.10%  "
$ws.Range("E42").Value = "  -0.45%  "
$ws.Range("E43").Value = "  +11.52%  "
$ws.Range("E44").Value = "  -5.76%  "
$ws.Range("E46").Value = "  -4.24%  "
$ws.Range("E47").Value = "  -1.54%  "
$ws.Range("E48").Value = "  -1.76%  "
$ws.Range("E49").Value = "  -1.02%  "
$ws.Range("E50").Value = "  -7.97%  "
$ws.Range("E51").Value = "  +6.97%  "
